{"js": "// Add a new sourced-link paragraph right after the \"Erstellung eines Main\n// Menus\" entry (and before the trailing empty paragraph), mirroring the\n// existing hyperlink-reference paragraphs already in the document:\n//   <tab-indent> <hyperlink>https://wallpaperaccess.com/full/38123.jpg</hyperlink> (Als Men\u00fc Hintergrund)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"Main Menu\" source entry so the new\n// paragraph lands right after it (and still before the trailing blank one).\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Erstellung eines Main Menus\") !== -1) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\nif (!anchor) {\n  throw new Error(\"Anchor paragraph 'Erstellung eines Main Menus' not found\");\n}\n\n// Insert a fresh empty paragraph right after the anchor.\nconst newPara = anchor.insertParagraph(\"\", \"After\");\n\n// Match the indentation used by the sibling source-link paragraphs\n// (<w:ind w:firstLine=\"708\"/> == 708 twips == 35.4 points).\nnewPara.paragraphFormat.firstLineIndent = 35.4;\n\nconst url = \"https://wallpaperaccess.com/full/38123.jpg\";\n\n// Insert the hyperlink's visible text first, then turn that range into an\n// actual hyperlink (applies the built-in \"Hyperlink\" character style).\nconst linkRange = newPara.insertText(url, \"Start\");\nlinkRange.hyperlink = url;\n\n// Trailing descriptive text after the link.\nnewPara.insertText(\" (Als Men\u00fc Hintergrund)\", \"End\");\n\nawait context.sync();\n", "ps1": "# Add a new sourced-link paragraph right after the \"Erstellung eines Main\n# Menus\" entry (and before the trailing empty paragraph), mirroring the\n# existing hyperlink-reference paragraphs already in the document:\n#   <tab-indent> <hyperlink>https://wallpaperaccess.com/full/38123.jpg</hyperlink> (Als Men\u00fc Hintergrund)\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph that ends the \"Main Menu\" source entry so the new\n# paragraph lands right after it (and still before the trailing blank one).\n$paras = $d.Paragraphs\n$anchor = $null\nfor ($i = 1; $i -le $paras.Count; $i++) {\n    $p = $paras.Item($i)\n    if ($p.Range.Text -like \"*Erstellung eines Main Menus*\") {\n        $anchor = $p\n        break\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph 'Erstellung eines Main Menus' not found\"\n}\n\n# Insert a fresh empty paragraph right after the anchor. It inherits the\n# anchor's paragraph formatting, i.e. the <w:ind w:firstLine=\"708\"/> used by\n# the other source-link paragraphs.\n$anchor.Range.InsertParagraphAfter()\n$newPara = $anchor.Next()\n\n$url = \"https://wallpaperaccess.com/full/38123.jpg\"\n\n# Write the hyperlink's visible text into the new (still empty) paragraph,\n# excluding its trailing paragraph mark, then convert that range into an\n# actual hyperlink (applies the built-in \"Hyperlink\" character style).\n$linkRange = $newPara.Range\n$linkRange.End = $linkRange.End - 1\n$linkRange.Text = $url\n$d.Hyperlinks.Add($linkRange, $url)\n\n# Append the trailing descriptive text right after the hyperlink, still\n# inside the same paragraph (before its paragraph mark).\n$tailRange = $newPara.Range\n$tailRange.End = $tailRange.End - 1\n$tailRange.Collapse(0)\n$tailRange.InsertAfter(\" (Als Men\u00fc Hintergrund)\")\n"}
